$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: add new "Start Date:" / date label pair in K1:L1 (mirrors B1:C1)
$ws.Range("K1").Value = "Start Date:"

# "07/14/2023" would be auto-parsed into a date serial if typed straight into
# a General-formatted cell, so force text formatting for the entry, then
# paste the original (General/style-2) formatting back over it so the cell's
# style matches its row-mates.
$ws.Range("L1").NumberFormat = "@"
$ws.Range("L1").Value = "07/14/2023"
$ws.Range("K1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 2: re-purpose J2:Q2 header labels
$ws.Range("J2").Value = "EQUIP_TYPE"
$ws.Range("K2").Value = "VENDOR"
$ws.Range("L2").Value = "UNIT_ID"
$ws.Range("N2").Value = "Start Date"
$ws.Range("O2").Value = "End Date"
$ws.Range("P2").Value = "Entity Name "
$ws.Range("Q2").Value = "Start Date "

# Selection moved from O6 to M6
$ws.Range("M6").Select()
